$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (2019/2020 columns revised) ---
$ws.Range("P5").Value = 25.6
$ws.Range("Q5").Value = 23.8

$ws.Range("P6").Value = 18.6
$ws.Range("Q6").Value = 16.7

$ws.Range("P8").Value = 2.1
$ws.Range("Q8").Value = 1.8

$ws.Range("P9").Value = 4.9
$ws.Range("Q9").Value = 5.2

# --- Add new columns R (2021) and S (2022) ---
# Row 4: year headers
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("Q4").Copy($ws.Range("S4"))
$ws.Range("R4").Value = 2021
$ws.Range("S4").Value = 2022

# Row 5: Revenues, total
$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("Q5").Copy($ws.Range("S5"))
$ws.Range("R5").Value = 26.8
$ws.Range("S5").Value = 26.8

# Row 6: Tax revenues
$ws.Range("Q6").Copy($ws.Range("R6"))
$ws.Range("Q6").Copy($ws.Range("S6"))
$ws.Range("R6").Value = 19.3
$ws.Range("S6").Value = 19.3

# Row 7: Received official transfers ("-")
$ws.Range("Q7").Copy($ws.Range("R7"))
$ws.Range("Q7").Copy($ws.Range("S7"))

# Row 8: Non-tax revenues
$ws.Range("Q8").Copy($ws.Range("R8"))
$ws.Range("Q8").Copy($ws.Range("S8"))
$ws.Range("R8").Value = 1.8
$ws.Range("S8").Value = 1.8

# Row 9: Revenues from the sale of non-financial assets
$ws.Range("Q9").Copy($ws.Range("R9"))
$ws.Range("Q9").Copy($ws.Range("S9"))
$ws.Range("R9").Value = 5.7
$ws.Range("S9").Value = 5.7

# Row 10: Contributions / deductions for social needs
$ws.Range("Q10").Copy($ws.Range("R10"))
$ws.Range("Q10").Copy($ws.Range("S10"))
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0

# --- Update selection to match new used range ---
$ws.Range("T3").Select()
